$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 47: backfill the Notes cell that was missing ---
$ws.Range("K47").Value = 'Lots of plutonium got stacked onto one big guy, but he failed to escape as quite a few villains are locations.'

# --- Row 48 ---
$ws.Range("A48").Value = 2
$ws.Range("B48").Value = 'Brainwash the Military'
$ws.Range("C48").Value = 'Dr. Strange'
$ws.Range("D48").Value = 'Defenders|Spider Friends'
$ws.Range("E48").Value = 'M.O.D.O.K.s'
$ws.Range("F48").Value = 'Steve Rogers, Director of S.H.I.E.L.D. (C75)|Venompool (VE)|Hulk (B)|Jessica Jones (D)|Deadpool (B)'
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = '18|45'
$ws.Range("I48").Value = 'no'
$ws.Range("K48").Value = 'No officer managed to escape. Steve Rogers and Venompool combine very strongly.'

# --- Row 49 ---
$ws.Range("A49").Value = 2
$ws.Range("B49").Value = 'Crush Them With My Bare Hands'
$ws.Range("C49").Value = 'Illuminati, Secret Society'
$ws.Range("D49").Value = 'Illuminati|Salvagers'
$ws.Range("E49").Value = 'Cape-Killers'
$ws.Range("F49").Value = 'Dr. Strange (SW1)|Totally Awesome Hulk (CH)|No-Name, Brood Queen (WW)|Korg (WW)|Phoenix (XM)'
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = '22|28'
$ws.Range("I49").Value = 'yes'
$ws.Range("K49").Value = 'Master strikes are really disruptive. The last ones came quite late, but we lost with one tactic remaining.'

# --- Row 50 ---
$ws.Range("A50").Value = 2
$ws.Range("B50").Value = 'S.H.I.E.L.D. vs. HYDRA War'
$ws.Range("C50").Value = 'General Ross'
$ws.Range("D50").Value = 'Hydra Elite|Code Red'
$ws.Range("E50").Value = 'Khonshu Guardians'
$ws.Range("F50").Value = 'Darkhawk (R)|Rogue (B)|Lady Sif (HOA)|War Machine (R)|Hercules (CW)'
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = '18|47'
$ws.Range("I50").Value = 'no'
$ws.Range("K50").Value = 'HYDRA level was 4 by ambush only.'

# --- Row 51 ---
$ws.Range("A51").Value = 2
$ws.Range("B51").Value = 'Fear Itself'
$ws.Range("C51").Value = 'Ultron'
$ws.Range("D51").Value = 'Ultron''s Legacy|Superhuman Registration Act'
$ws.Range("E51").Value = 'Cytoplasm Spikes'
$ws.Range("F51").Value = 'Namora (WW)|Hercules (CW)|Lady Thor (SW1)|Phoenix (XM)|Beta Ray Bill (HOA)'
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = '26|45'
$ws.Range("I51").Value = 'no'
$ws.Range("K51").Value = 'Fear level got to 3. Phoenix kept the city empty, Lady Thor netted a lot of recruit to buy ultimates. Master strikes only ensured empowered to green.'

# --- Row 52 ---
$ws.Range("A52").Value = 2
$ws.Range("B52").Value = 'Secret Wars'
$ws.Range("C52").Value = 'Mole Man'
$ws.Range("D52").Value = 'Subterranea|Sinister Six'
$ws.Range("E52").Value = 'Khonshu Guardians'
$ws.Range("F52").Value = 'Ruby Summers (SW2)|Sentry (WW)|Skirn, Breaker of Men (FI)|Jubilee (XM)|Ant-Man (AM)'
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = '42|27'
$ws.Range("I52").Value = 'yes'
$ws.Range("J52").Value = 'Random masterminds were Magneto, Arcade and The Red King. Arcade got his human shields when entering.'
$ws.Range("K52").Value = 'Only one twist out, but it got put close to the bottom by one of Mole man''s tactics.'

# --- Row 53 ---
$ws.Range("A53").Value = 2
$ws.Range("B53").Value = 'The Korvac Saga'
$ws.Range("C53").Value = 'The Hood'
$ws.Range("D53").Value = 'Hood''s Gang|Code Red'
$ws.Range("E53").Value = 'Doombot Legion'
$ws.Range("F53").Value = 'Bob, Agent of HYDRA (DP)|Mockingbird (AOS)|Goliath (CW)|Wolfsbane (NM)|Spider-Gwen (SW2)'
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = '12|42'
$ws.Range("I53").Value = 'no'
$ws.Range("K53").Value = 'Early spider-gwen is insane.'

# --- Row 54 ---
$ws.Range("A54").Value = 2
$ws.Range("B54").Value = 'Mutant-Hunting Super Sentinels'
$ws.Range("C54").Value = 'Stryfe'
$ws.Range("D54").Value = 'MLF|Infinity Gems'
$ws.Range("E54").Value = 'Sentinel|Death''s Heads'
$ws.Range("F54").Value = 'She-Hulk (WW)|Luke Cage (CW)|Wolverine (B)|Iron Fist (DC)|X-23 (XM)'
$ws.Range("G54").Value = 1
$ws.Range("H54").Value = '7|35'
$ws.Range("I54").Value = 'no'
$ws.Range("K54").Value = 'Iron Fist, She-Hulk and Wolverine are insane together. Stryfe got pummeled as one player got up to 59 attack. No Sentinel escaped.'

# --- Row 55 ---
$ws.Range("A55").Value = 2
$ws.Range("B55").Value = 'Smash Two Dimensions Together'
$ws.Range("C55").Value = 'The Goblin, Underworld Boss'
$ws.Range("D55").Value = 'Goblin''s Freak Show|K''un-lun|Life Foundation'
$ws.Range("E55").Value = 'Magma Men'
$ws.Range("F55").Value = 'Totally Awesome Hulk (CH)|Wasp (AM)|Gambit (B)|Spider-Man Noir (N)|Venompool (VE)'
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = '75|43'
$ws.Range("I55").Value = 'no'
$ws.Range("K55").Value = 'Wasp and Gambit combo. Also lots of bystanders from Spiderman Noir and Venompool. Only two escapes, but Goblin took a long time to defeat.'

# --- Row 56 ---
$ws.Range("A56").Value = 2
$ws.Range("B56").Value = 'Predict Future Crime'
$ws.Range("C56").Value = 'Ultron'
$ws.Range("D56").Value = 'Ultron''s Legacy|Shadow-X|Warbound'
$ws.Range("E56").Value = 'Mandarin''s Rings'
$ws.Range("F56").Value = 'Star-Lord (GG)|Agent Phil Coulson (AOS)|Juggernaut (V)|Psylocke (XM)|Legion (XM)'
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = '7|8'
$ws.Range("I56").Value = 'yes'
$ws.Range("K56").Value = 'Three villains escaped because almost all twists passed in the first ten reveals, with no henchmen revealed. A late twist hit three villains and it was game over.'

# --- Row 57 ---
$ws.Range("A57").Value = 2
$ws.Range("B57").Value = 'The Demon Bear Saga'
$ws.Range("C57").Value = 'Illuminati, Secret Society'
$ws.Range("D57").Value = 'Demons of Limbo|Illuminati'
$ws.Range("E57").Value = 'Savage Land Mutates'
$ws.Range("F57").Value = 'Cannonball (XM)|Beast (SW2)|Deadpool (DP)|Captain America, Secret Avenger (CW)|Hulk (B)'
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = '29|44'
$ws.Range("I57").Value = 'yes'
$ws.Range("K57").Value = 'No KO mechanisms made the game very slow. Early master strikes really screwed the early recruits. Game ended because villain deck ran out (bear escaped twice).'

# --- Update the view state: scroll position + active selection ---
[void]$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 31
$win.ScrollColumn = 5
[void]$ws.Range("H58").Select()
